$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.356322765350342
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 0.7371582388877869
$ws.Range("E1").Value = 0.902252197265625
